$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Hành động của sinh viên..." paragraph: the sentence was split across
#    two runs around the _GoBack bookmark ("...đi học, là" | bookmark |
#    "m bài tập, nộp học phí..."). Merge them into a single run with the
#    full sentence, which also removes the old (collapsed) bookmark that
#    sat between the two runs.
# ---------------------------------------------------------------------------
$rStart = $d.Content.Duplicate
$rStart.Find.Execute("Hành động của sinh viên gồm: đi học, là") | Out-Null
$mergeBegin = $rStart.Start

$rEnd = $d.Content.Duplicate
$rEnd.Find.Execute("m bài tập, nộp học phí, làm bài kiểm tra và cho biết thông tin đầy đủ của bản thân, các hàm tạo có và không có tham số, các phương thức get/set.") | Out-Null
$mergeEnd = $rEnd.End

$mergeRange = $d.Range($mergeBegin, $mergeEnd)
$mergeRange.Text = "Hành động của sinh viên gồm: đi học, làm bài tập, nộp học phí, làm bài kiểm tra và cho biết thông tin đầy đủ của bản thân, các hàm tạo có và không có tham số, các phương thức get/set."

# ---------------------------------------------------------------------------
# 2. Insert a space between ">=" and "1" in the credit-count condition.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
  "số tín chỉ(>=1 và",
  $true, $false, $false, $false, $false,
  $true, 1, $false,
  "số tín chỉ(>= 1 và",
  2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Re-create the _GoBack bookmark at the end of the "Ex2:" paragraph
#    (right after its last run, before the paragraph mark). A collapsed
#    bookmark placed directly at that end-of-paragraph-content offset gets
#    mis-resolved by this host, so we build it around a throw-away
#    placeholder character and then delete that character - the bookmark
#    then correctly collapses back down to a zero-length bookmark.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}

$ex2End = $d.Content.Duplicate
$ex2End.Find.Execute(", người dạy. Môn học có thể thực hiện các hành động như cung cấp đầy đủ thông tin bản thân, có các constructor, getter/setter tương ứng các thuộc tính.") | Out-Null
$insertPos = $ex2End.End

$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanup = $d.Range($insertPos, $insertPos + 1)
$cleanup.Text = ""
